$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2 ("User Information") - rebuild the reference-data block
# ---------------------------------------------------------------------------

# Turn the existing table (Table2) back into a normal range, keeping data
$ws2.ListObjects.Item(1).Unlist()

# Remove the old "postal-code / 12345" row entirely (shifts rows up)
$ws2.Rows.Item(5).Delete()

# Clear the old labels/values that are being replaced
$ws2.Range("A2:B4").ClearContents()

# New header row (typed right-to-left, matching the resulting shared-string order)
$ws2.Range("C1").Value = "Postalcode"
$ws2.Range("B1").Value = "Lastname"
$ws2.Range("A1").Value = "Firstname"

# New data rows
$ws2.Range("A2").Value = "John"
$ws2.Range("B2").Value = "Test"
$ws2.Range("C2").Value = "JHRWG"

$ws2.Range("A3").Value = "Sarah"
$ws2.Range("B3").Value = "XXXTest"
$ws2.Range("C3").Value = "ABCDEF"

# Re-apply the left-aligned style that the shifted empty cell (B4) should carry
$ws2.Range("B4").HorizontalAlignment = -4131

# Header row fill (theme accent2 / theme index 5)
$ws2.Range("A1:C1").Interior.ThemeColor = 6

# Data rows fill (theme accent6 / theme index 9)
$ws2.Range("A2:B3").Interior.ThemeColor = 10
$ws2.Range("C2").Interior.ThemeColor = 10
$ws2.Range("C3").Interior.ThemeColor = 10

# The postal-code value cell (C2) is also left aligned
$ws2.Range("C2").HorizontalAlignment = -4131

# Column widths / layout tweaks
$ws2.Columns.Item(1).ColumnWidth = 19.59244791666667
$ws2.Columns.Item(3).ColumnWidth = 12.451822916666666

# Selection / view bookkeeping
$ws2.Range("E7").Select()

# ---------------------------------------------------------------------------
# Sheet1 ("Login Details") - selection bookkeeping only
# ---------------------------------------------------------------------------
$ws1.Range("B12").Select()

Write-Host "done"
